$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: change its formatting group to the "bottom border" variant (styles 8/8/9/9/9),
#     matching the look used for row 6 (the other group-separator row). Values are unchanged.
$ws.Range("A6:E6").Copy() | Out-Null
$ws.Range("A8:E8").PasteSpecial(-4122) | Out-Null

# --- Row 9: new data row, reuses row 7's formatting (styles 4/4/5/5/5).
$ws.Range("A7:E7").Copy() | Out-Null
$ws.Range("A9:E9").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(9).RowHeight = 43.2

$ws.Range("A9").Value = "SCRIPT/T01P02A/um1318.ssb"
$ws.Range("B9").Value = 89
$ws.Range("C9").Value = " They say if you take a Time\nGear away, the time stops where that Time\nGear was."
$ws.Range("D9").Value = " Говорят, что если убрать\nШестерню Времени, то там, где она была,\nостановится время."
$ws.Range("E9").Value = " Ãïâïñÿó, œóï åòìé ôáñàóû\nŠåòóåñîý Âñåíåîé, óï óàí, ãäå ïîà áúìà,\nïòóàîïâéóòÿ âñåíÿ."

# --- Row 10: only column A is filled in (same style family as row 9 / row 7).
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(10).RowHeight = 43.2

$ws.Range("A10").Value = "SCRIPT/T01P02A/um1404.ssb"

$excel.CutCopyMode = $false

# --- Update the view to reflect the final scroll/selection position.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("C9").Select() | Out-Null

Write-Output "edit applied"
